$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the manualStatus values in column I from plain numbers to a
# bracketed text status, e.g. 32 -> "[32]", 4 -> "[4]"
$ws.Range("I5").Value = "[32]"
$ws.Range("I12").Value = "[4]"
$ws.Range("I14").Value = "[4]"

# Widen column F (fastqFileName) to fit the long file names
$ws.Columns.Item(6).ColumnWidth = 58.22

# Move the active selection to I15
$ws.Range("I15").Select()
